$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header
# formatting (bold, centered, bordered) used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I ("I0") and J ("IF") for rows 2-7.
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 12

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 4

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 7
